# Apply cryptos list update (Tue Feb 27 08:42:06 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.390.62'
$ws.Range("E2").Value = '  +9.86%  '
$ws.Range("D3").Value = '3.227.79'
$ws.Range("E3").Value = '  +4.14%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '397.79'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.18'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +7.10%  '
$ws.Range("E7").Value = '  +2.84%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.619'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +5.99%  '
$ws.Range("E10").Value = '  +6.59%  '
$ws.Range("E11").Value = '  +6.82%  '
$ws.Range("E12").Value = '  +2.09%  '
$ws.Range("D13").Value = '3.740.77'
$ws.Range("E13").Value = '  +4.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.10'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.06'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.87%  '
$ws.Range("D16").Value = '3.219.50'
$ws.Range("E16").Value = '  +3.65%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.04'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +4.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.90'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.65%  '
$ws.Range("D19").Value = '56.203.06'
$ws.Range("E19").Value = '  +9.26%  '
$ws.Range("E20").Value = '  +3.34%  '
$ws.Range("E21").Value = '  +6.44%  '
$ws.Range("E22").Value = '  +4.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '298.46'
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.70'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +8.15%  '
$ws.Range("E25").Value = '  +1.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.15'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.08'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.44'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.90%  '
$ws.Range("E29").Value = '  +4.57%  '
$ws.Range("E30").Value = '  +0.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.110'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.14'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +6.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0495'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '36.72'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.13'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.45%  '
$ws.Range("E36").Value = '  +3.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.13'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +25.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.53'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +4.40%  '
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '134.84'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.47'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +4.74%  '
$ws.Range("E42").Value = '  +3.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.00'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +3.36%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.119'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.08%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.285'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.23'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +55.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.26'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.78%  '
$ws.Range("E48").Value = '  -2.17%  '
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("D50").Value = '2.131.75'
$ws.Range("E50").Value = '  +2.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0362'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +10.63%  '
